$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the people in column A (shared-string text swap)
$ws.Range("A2").Value = "Assunta Feil"
$ws.Range("A3").Value = "Al Greenholt"
$ws.Range("A4").Value = "Keven Mueller"
$ws.Range("A5").Value = "Katharina Shields"

# Fix the popularity-bar values: row 3 (Al Greenholt) is no longer
# applied/accepted, row 5 (Katharina Shields) now is.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Column A widened to fit the new, longer name ("Katharina Shields")
$ws.Columns.Item(1).ColumnWidth = 20.25
